$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 through 22 had their custom height reduced from 19.5 to 18.75
$ws.Range("A2:A22").RowHeight = 18.75
